$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.733.11'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '2.086.84'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'232.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = "'0.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = "'57.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.32%  '
$ws.Range("E9").Value = '  +1.57%  '
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("E11").Value = '  +3.00%  '
$ws.Range("D12").Value = '2.381.19'
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").Value = "'21.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.02%  '
$ws.Range("D15").Value = "'0.767"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("E16").Value = '  +1.95%  '
$ws.Range("D17").Value = '2.083.64'
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").Value = '37.615.95'
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("D19").Value = "'6.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.05%  '
$ws.Range("D20").Value = "'70.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.15%  '
$ws.Range("D21").Value = '0.0₃0820'
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("D22").Value = "'227.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = "'2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("D26").Value = "'168.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("E27").Value = '  +10.03%  '
$ws.Range("D28").Value = "'8.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.78%  '
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("D30").Value = "'19.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.22%  '
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("D32").Value = "'4.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.30%  '
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("D34").Value = "'4.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.34%  '
$ws.Range("D35").Value = "'2.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("E36").Value = '  +4.23%  '
$ws.Range("D37").Value = "'3.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.60%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -5.21%  '
$ws.Range("D40").Value = "'0.0995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.77%  '
$ws.Range("E41").Value = '  -0.31%  '
$ws.Range("D42").Value = "'96.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("E43").Value = '  +0.36%  '
$ws.Range("D44").Value = '1.451.28'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("E46").Value = '  +3.44%  '
$ws.Range("D47").Value = "'4.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.50%  '
$ws.Range("D48").Value = "'15.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.63%  '
$ws.Range("D49").Value = "'7.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.27%  '
$ws.Range("D50").Value = "'3.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.94%  '
$ws.Range("D51").Value = '2.276.66'
$ws.Range("E51").Value = '  +1.14%  '
